$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '47.392.15'
$ws.Range("E2").Value = '  +2.97%  '

# Row 3
$ws.Range("D3").Value = '2.510.56'
$ws.Range("E3").Value = '  +2.59%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = '''110.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.86%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '''324.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.98%  '

# Row 7
$ws.Range("E7").Value = '  +1.59%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("E9").Value = '  +1.50%  '

# Row 10
$ws.Range("D10").Value = '''39.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.68%  '

# Row 11
$ws.Range("E11").Value = '  +1.81%  '

# Row 12
$ws.Range("E12").Value = '  +0.99%  '

# Row 13
$ws.Range("E13").Value = '  +1.89%  '

# Row 14
$ws.Range("D14").Value = '''7.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.76%  '

# Row 15
$ws.Range("D15").Value = '2.903.77'
$ws.Range("E15").Value = '  +2.57%  '

# Row 16
$ws.Range("D16").Value = '2.511.90'
$ws.Range("E16").Value = '  +1.64%  '

# Row 17
$ws.Range("D17").Value = '''0.862'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.77%  '

# Row 18
$ws.Range("D18").Value = '47.359.50'
$ws.Range("E18").Value = '  +3.16%  '

# Row 19
$ws.Range("D19").Value = '''12.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.57%  '

# Row 20
$ws.Range("D20").Value = '''6.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.76%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0949'
$ws.Range("E21").Value = '  +1.58%  '

# Row 22
$ws.Range("E22").Value = '  +13.03%  '

# Row 23
$ws.Range("D23").Value = '''70.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.63%  '

# Row 24
$ws.Range("D24").Value = '''250.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.26%  '

# Row 25
$ws.Range("E25").Value = '  +4.33%  '

# Row 26
$ws.Range("D26").Value = '''26.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.89%  '

# Row 27
$ws.Range("D27").Value = '''0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.06%  '

# Row 28
$ws.Range("D28").Value = '''2.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.61%  '

# Row 29
$ws.Range("D29").Value = '''10.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.58%  '

# Row 30
$ws.Range("D30").Value = '''35.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.83%  '

# Row 31
$ws.Range("D31").Value = '''0.137'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.12%  '

# Row 32
$ws.Range("D32").Value = '''50.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.88%  '

# Row 33
$ws.Range("D33").Value = '''19.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.19%  '

# Row 34
$ws.Range("D34").Value = '''5.46'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.63%  '

# Row 35
$ws.Range("D35").Value = '''0.0799'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.91%  '

# Row 36
$ws.Range("E36").Value = '  +0.16%  '

# Row 37
$ws.Range("D37").Value = '''2.01'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.99%  '

# Row 38
$ws.Range("D38").Value = '''4.74'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.35%  '

# Row 39
$ws.Range("D39").Value = '''3.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.91%  '

# Row 40
$ws.Range("E40").Value = '  +1.69%  '

# Row 41
$ws.Range("D41").Value = '''122.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.81%  '

# Row 42
$ws.Range("E42").Value = '  -0.96%  '

# Row 43
$ws.Range("D43").Value = '''21.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.97%  '

# Row 44
$ws.Range("E44").Value = '  +2.51%  '

# Row 45
$ws.Range("D45").Value = '2.002.73'
$ws.Range("E45").Value = '  +2.22%  '

# Row 46
$ws.Range("E46").Value = '  +5.15%  '

# Row 47
$ws.Range("E47").Value = '  -1.40%  '

# Row 48
$ws.Range("D48").Value = '''1.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.12%  '

# Row 49
$ws.Range("D49").Value = '''9.08'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.25%  '

# Row 50
$ws.Range("D50").Value = '''5.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.94%  '

# Row 51
$ws.Range("D51").Value = '''78.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
